$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header rich-text strings: bulletin volume/number and week date-range ---
# (Diff only retargets the trailing run of each shared string; the runtime
#  flattens rich-text runs on write, but every run in these two strings shares
#  identical formatting, so writing the full string is visually equivalent.)
$ws.Range("A8").Value = "Volume 30   Number  19"
$ws.Range("C9").Value = "Report Covering the Week  5/8/2023  Through  5/14/2023"

# --- Weekly crime-stat table (rows 15-30): refreshed counts/percentages ---

# Row 15
$ws.Range("L15").Value = -50
$ws.Range("N15").Value = -50

# Row 16
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 50
$ws.Range("F16").Value = 12
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = 9.090909090909
$ws.Range("I16").Value = 44
$ws.Range("J16").Value = 49
$ws.Range("K16").Value = -10.204081632653
$ws.Range("L16").Value = -13.725490196078
$ws.Range("M16").Value = -42.105263157894
$ws.Range("N16").Value = -89.573459715639

# Row 17
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 2
$ws.Range("F17").Value = 15
$ws.Range("G17").Value = 13
$ws.Range("H17").Value = 15.384615384615
$ws.Range("I17").Value = 90
$ws.Range("J17").Value = 57
$ws.Range("K17").Value = 57.894736842105
$ws.Range("L17").Value = 83.673469387755
$ws.Range("M17").Value = 91.489361702127
$ws.Range("N17").Value = -38.775510204081

# Row 18
$ws.Range("C18").Value = 6
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 200
$ws.Range("F18").Value = 16
$ws.Range("H18").Value = 23.076923076923
$ws.Range("I18").Value = 84
$ws.Range("J18").Value = 77
$ws.Range("K18").Value = 9.090909090909
$ws.Range("L18").Value = 86.666666666666
$ws.Range("M18").Value = 140
$ws.Range("N18").Value = -68.773234200743

# Row 19
$ws.Range("C19").Value = 16
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = 60
$ws.Range("F19").Value = 59
$ws.Range("G19").Value = 46
$ws.Range("H19").Value = 28.260869565217
$ws.Range("I19").Value = 248
$ws.Range("J19").Value = 214
$ws.Range("K19").Value = 15.887850467289
$ws.Range("L19").Value = 43.352601156069
$ws.Range("M19").Value = 61.038961038961
$ws.Range("N19").Value = -30.140845070422

# Row 20
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 50
$ws.Range("F20").Value = 10
$ws.Range("G20").Value = 4
$ws.Range("H20").Value = 150
$ws.Range("I20").Value = 26
$ws.Range("J20").Value = 19
$ws.Range("K20").Value = 36.842105263157
$ws.Range("L20").Value = 225
$ws.Range("M20").Value = 44.444444444444
$ws.Range("N20").Value = -89.723320158102

# Row 21
$ws.Range("C21").Value = 30
$ws.Range("D21").Value = 18
$ws.Range("E21").Value = 66.666666666666
$ws.Range("F21").Value = 112
$ws.Range("G21").Value = 88
$ws.Range("H21").Value = 27.272727272727
$ws.Range("I21").Value = 496
$ws.Range("J21").Value = 421
$ws.Range("K21").Value = 17.814726840855
$ws.Range("L21").Value = 49.848942598187
$ws.Range("M21").Value = 49.848942598187
$ws.Range("N21").Value = -65.887207702888

# Row 22
$ws.Range("F22").Value = 7
$ws.Range("H22").Value = 75
$ws.Range("I22").Value = 21
$ws.Range("J22").Value = 16
$ws.Range("K22").Value = 31.25
$ws.Range("L22").Value = 5
$ws.Range("M22").Value = -19.230769230769
$ws.Range("D22").NumberFormat = '#,##0'
$ws.Range("D22").Value = 1
$ws.Range("E22").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E22").Value = 200

# Row 23
$ws.Range("F23").Value = 6
$ws.Range("G23").Value = 3
$ws.Range("H23").Value = 100
$ws.Range("I23").Value = 20
$ws.Range("K23").Value = 25
$ws.Range("L23").Value = 11.111111111111
$ws.Range("M23").Value = 233.333333333333
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "0"
$ws.Range("C23").NumberFormat = "general"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0"
$ws.Range("D23").NumberFormat = "general"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "***.*"
$ws.Range("E23").NumberFormat = "general"

# Row 24
$ws.Range("C24").Value = 39
$ws.Range("D24").Value = 36
$ws.Range("E24").Value = 8.333333333333
$ws.Range("F24").Value = 179
$ws.Range("G24").Value = 140
$ws.Range("H24").Value = 27.857142857142
$ws.Range("I24").Value = 733
$ws.Range("J24").Value = 596
$ws.Range("K24").Value = 22.986577181208
$ws.Range("L24").Value = 46.89378757515
$ws.Range("M24").Value = 50.823045267489

# Row 25
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 14
$ws.Range("E25").Value = -78.571428571428
$ws.Range("F25").Value = 26
$ws.Range("G25").Value = 39
$ws.Range("H25").Value = -33.333333333333
$ws.Range("I25").Value = 136
$ws.Range("J25").Value = 133
$ws.Range("K25").Value = 2.255639097744
$ws.Range("L25").Value = 60
$ws.Range("M25").Value = -6.849315068493

# Row 26
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = -50
$ws.Range("I26").Value = 3
$ws.Range("J26").Value = 6
$ws.Range("K26").Value = -50
$ws.Range("L26").Value = -57.142857142857
$ws.Range("C26").NumberFormat = '#,##0'
$ws.Range("C26").Value = 1
$ws.Range("D26").NumberFormat = '#,##0'
$ws.Range("D26").Value = 1
$ws.Range("E26").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E26").Value = 0
$ws.Range("F26").NumberFormat = '#,##0'
$ws.Range("F26").Value = 1

# Row 27
$ws.Range("E27").Value = -100
$ws.Range("J27").Value = 26
$ws.Range("K27").Value = -34.615384615384
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "0"
$ws.Range("C27").NumberFormat = "general"

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0"
$ws.Range("D28").NumberFormat = "general"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "***.*"
$ws.Range("E28").NumberFormat = "general"

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0"
$ws.Range("D29").NumberFormat = "general"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "***.*"
$ws.Range("E29").NumberFormat = "general"

# Row 30
$ws.Range("L30").Value = 233.333333333333
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "0"
$ws.Range("C30").NumberFormat = "general"
